# Auto-generated edit script: updates cryptos list per Fri Dec  1 12:45:47 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Val)
    # Force the cell to stay text even when the new value parses as a number,
    # mirroring how these price cells were authored (inline strings), then drop
    # back to the default style so no stray number-format/quote-prefix is left behind.
    $Sheet.Range($Addr).NumberFormat = "@"
    $Sheet.Range($Addr).Value = $Val
    $Sheet.Range($Addr).Style = "Normal"
}

$ws.Range('D2').Value = '38.488.83'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '2.090.71'
$ws.Range('E3').Value = '  +2.14%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue $ws 'D5' '228.41'
$ws.Range('E5').Value = '  +0.10%  '
Set-TextValue $ws 'D6' '0.612'
$ws.Range('E6').Value = '  +0.59%  '
Set-TextValue $ws 'D7' '61.01'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.10%  '
Set-TextValue $ws 'D9' '0.381'
$ws.Range('E9').Value = '  +1.17%  '
Set-TextValue $ws 'D10' '0.0838'
$ws.Range('E10').Value = '  +2.30%  '
Set-TextValue $ws 'D11' '0.104'
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').Value = '2.401.28'
$ws.Range('E12').Value = '  +2.15%  '
Set-TextValue $ws 'D13' '14.85'
$ws.Range('E13').Value = '  +1.32%  '
Set-TextValue $ws 'D14' '22.35'
$ws.Range('E14').Value = '  +6.14%  '
Set-TextValue $ws 'D15' '0.784'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('E16').Value = '  +4.68%  '
$ws.Range('D17').Value = '2.091.72'
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('D18').Value = '38.345.58'
$ws.Range('E18').Value = '  +1.28%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws 'D19' '6.07'
$ws.Range('E19').Value = '  +2.75%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws 'D20' '70.94'
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('E21').Value = '  +1.23%  '
Set-TextValue $ws 'D22' '225.73'
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('E23').Value = '  +0.00%  '
Set-TextValue $ws 'D24' '2.44'
$ws.Range('E24').Value = '  -0.08%  '
Set-TextValue $ws 'D25' '2.31'
$ws.Range('E25').Value = '  +1.38%  '
Set-TextValue $ws 'D26' '169.51'
$ws.Range('E26').Value = '  +0.93%  '
Set-TextValue $ws 'D27' '9.44'
$ws.Range('E27').Value = '  +0.99%  '
Set-TextValue $ws 'D28' '0.135'
$ws.Range('E28').Value = '  +4.48%  '
Set-TextValue $ws 'D29' '19.05'
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('E30').Value = '  +6.43%  '
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('E32').Value = '  +6.13%  '
Set-TextValue $ws 'D33' '4.77'
$ws.Range('E33').Value = '  +5.33%  '
Set-TextValue $ws 'D34' '4.53'
$ws.Range('E34').Value = '  +3.04%  '
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws 'D36' '2.39'
$ws.Range('E36').Value = '  +2.20%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws 'D37' '6.42'
$ws.Range('E37').Value = '  -2.73%  '
Set-TextValue $ws 'D38' '3.58'
$ws.Range('E38').Value = '  +4.02%  '
$ws.Range('E39').Value = '  -0.11%  '
Set-TextValue $ws 'D40' '18.35'
$ws.Range('E40').Value = '  +1.67%  '
$ws.Range('D41').Value = '1.540.49'
$ws.Range('E41').Value = '  -0.30%  '
Set-TextValue $ws 'D42' '99.85'
$ws.Range('E42').Value = '  +3.63%  '
$ws.Range('E43').Value = '  +1.40%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D44' '0.0939'
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D45' '2.83'
$ws.Range('E45').Value = '  +0.79%  '
Set-TextValue $ws 'D46' '7.84'
$ws.Range('E46').Value = '  +11.48%  '
Set-TextValue $ws 'D47' '4.18'
$ws.Range('E47').Value = '  +1.92%  '
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('E50').Value = '  +0.71%  '
$ws.Range('D51').Value = '2.287.48'
$ws.Range('E51').Value = '  +2.19%  '
